# Update cryptos list (prices, volumes, and a couple of reordered rows)
# with the latest data, as scraped by GitHub Actions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.886.15"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.95%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.268.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.85%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.65%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.91%  "

# Row 7
$ws.Range("E7").Value = "  -0.06%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.604"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.268.12"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.88%  "

# Row 10
$ws.Range("E10").Value = "  +7.34%  "

# Row 11
$ws.Range("E11").Value = "  +1.99%  "

# Row 12
$ws.Range("E12").Value = "  +6.07%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.835.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.78%  "

# Row 14
$ws.Range("E14").Value = "  +1.21%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.93%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.871.47"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.97%  "

# Row 17
$ws.Range("E17").Value = "  +3.21%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.260.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.54%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.03%  "

# Row 20
$ws.Range("E20").Value = "  +4.72%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.93%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.45%  "

# Row 23
$ws.Range("E23").Value = "  -0.08%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.15"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.49%  "

# Row 25
$ws.Range("E25").Value = "  +3.70%  "

# Row 26
$ws.Range("E26").Value = "  +4.91%  "

# Row 27
$ws.Range("E27").Value = "  -2.30%  "

# Row 28
$ws.Range("E28").Value = "  +2.02%  "

# Row 29
$ws.Range("E29").Value = "  +0.02%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.99"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.24%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.40%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.34%  "

# Row 33
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$ws.Range("E34").Value = "  +6.06%  "

# Row 35
$ws.Range("E35").Value = "  +4.65%  "

# Row 36
$ws.Range("E36").Value = "  +4.94%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.58"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.88%  "

# Row 38
$ws.Range("E38").Value = "  +1.63%  "

# Row 39
$ws.Range("E39").Value = "  +2.51%  "

# Row 40: RenderToken -> EnergySwap (name/link/price/volume swap)
$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.90"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.64%  "

# Row 41: EnergySwap -> RenderToken (name/link/price/volume swap)
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.16%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.46%  "

# Row 43
$ws.Range("E43").Value = "  +5.40%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "351.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.78%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.679.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.88%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.65"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.72%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.85"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.87%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0682"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.43%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0282"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.51%  "

# Row 50
$ws.Range("E50").Value = "  +5.38%  "

# Row 51
$ws.Range("E51").Value = "  +0.88%  "
